# The presentation ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  ("Office Theme" colours) - used by the Notes Master
#   ppt/theme/theme2.xml  ("Integral" colours)      - used by the Slide Master / slides
#
# The authored edit swaps the two themes' content wholesale (file names /
# relationships stay put, only the payload moves): theme1.xml ends up with
# the "Integral" colours and theme2.xml ends up with the "Office Theme"
# colours. Both themes already share an identical font scheme and format
# scheme, so the only substantive difference is the 12-colour colour
# scheme (clrScheme) carried by each part.
#
# The PowerPoint object model exposes that colour scheme as the 12-entry
# ThemeColorScheme on a Slide/SlideRange (it writes straight through to the
# slide master's theme part, i.e. theme2.xml here). Reassign every entry so
# theme2.xml ends up holding the colours that theme1.xml used to have
# ("Office" palette).

function Set-ThemeColor($scheme, $index, $r, $g, $b) {
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$tcs = $slide.ThemeColorScheme

# Index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
Set-ThemeColor $tcs 1  0x00 0x00 0x00   # dk1      (unchanged)
Set-ThemeColor $tcs 2  0xFF 0xFF 0xFF   # lt1      (unchanged)
Set-ThemeColor $tcs 3  0x44 0x54 0x6A   # dk2      -> 44546A
Set-ThemeColor $tcs 4  0xE7 0xE6 0xE6   # lt2      -> E7E6E6
Set-ThemeColor $tcs 5  0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
Set-ThemeColor $tcs 6  0xED 0x7D 0x31   # accent2  -> ED7D31
Set-ThemeColor $tcs 7  0xA5 0xA5 0xA5   # accent3  -> A5A5A5
Set-ThemeColor $tcs 8  0xFF 0xC0 0x00   # accent4  -> FFC000
Set-ThemeColor $tcs 9  0x44 0x72 0xC4   # accent5  -> 4472C4
Set-ThemeColor $tcs 10 0x70 0xAD 0x47   # accent6  -> 70AD47
Set-ThemeColor $tcs 11 0x05 0x63 0xC1   # hlink    -> 0563C1
Set-ThemeColor $tcs 12 0x95 0x4F 0x72   # folHlink -> 954F72
